$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell carrying the plain/default (unstyled) format used by all
# B:E data cells in this sheet - used to reset style after forcing a
# numeric-looking value to be stored as text (NumberFormat = "@" otherwise
# leaves a residual Text number-format on the cell).
$normalStyleTemplate = $ws.Range("B2")

$ws.Range('D2').Value = '68.020.53'
$ws.Range('E2').Value = '  +1.81%  '

$ws.Range('D3').Value = '3.333.05'
$ws.Range('E3').Value = '  +1.52%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $normalStyleTemplate.Style
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.58'
$ws.Range('D5').Style = $normalStyleTemplate.Style
$ws.Range('E5').Value = '  +1.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.63'
$ws.Range('D6').Style = $normalStyleTemplate.Style
$ws.Range('E6').Value = '  +1.48%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('E8').Value = '  +1.61%  '

$ws.Range('D9').Value = '3.330.57'
$ws.Range('E9').Value = '  +1.53%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.185'
$ws.Range('D10').Style = $normalStyleTemplate.Style
$ws.Range('E10').Value = '  +6.31%  '

$ws.Range('E11').Value = '  +1.67%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.25'
$ws.Range('D12').Style = $normalStyleTemplate.Style
$ws.Range('E12').Value = '  +3.94%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('D13').Style = $normalStyleTemplate.Style
$ws.Range('E13').Value = '  +2.50%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '691.44'
$ws.Range('D14').Style = $normalStyleTemplate.Style
$ws.Range('E14').Value = '  +0.40%  '

$ws.Range('D15').Value = '3.868.92'
$ws.Range('E15').Value = '  +1.53%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.43'
$ws.Range('D16').Style = $normalStyleTemplate.Style
$ws.Range('E16').Value = '  +1.83%  '

$ws.Range('D17').Value = '67.952.41'
$ws.Range('E17').Value = '  +1.56%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.118'
$ws.Range('D18').Style = $normalStyleTemplate.Style
$ws.Range('E18').Value = '  -0.45%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.335.57'
$ws.Range('E19').Value = '  +1.55%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.47'
$ws.Range('D20').Style = $normalStyleTemplate.Style
$ws.Range('E20').Value = '  +1.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.09'
$ws.Range('D21').Style = $normalStyleTemplate.Style
$ws.Range('E21').Value = '  +3.53%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.898'
$ws.Range('D22').Style = $normalStyleTemplate.Style
$ws.Range('E22').Value = '  +1.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.41'
$ws.Range('D23').Style = $normalStyleTemplate.Style
$ws.Range('E23').Value = '  +4.77%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.06'
$ws.Range('D24').Style = $normalStyleTemplate.Style
$ws.Range('E24').Value = '  +0.50%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.62'
$ws.Range('D25').Style = $normalStyleTemplate.Style
$ws.Range('E25').Value = '  +0.79%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.92'
$ws.Range('D26').Style = $normalStyleTemplate.Style
$ws.Range('E26').Value = '  +1.24%  '

$ws.Range('E27').Value = '  +0.20%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.59'
$ws.Range('D28').Style = $normalStyleTemplate.Style
$ws.Range('E28').Value = '  +3.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.24'
$ws.Range('D29').Style = $normalStyleTemplate.Style
$ws.Range('E29').Value = '  -0.90%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.58'
$ws.Range('D30').Style = $normalStyleTemplate.Style
$ws.Range('E30').Value = '  +2.64%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.12'
$ws.Range('D31').Style = $normalStyleTemplate.Style
$ws.Range('E31').Value = '  +5.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '567.49'
$ws.Range('D32').Style = $normalStyleTemplate.Style
$ws.Range('E32').Value = '  -0.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.02'
$ws.Range('D33').Style = $normalStyleTemplate.Style
$ws.Range('E33').Value = '  +1.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.106'
$ws.Range('D34').Style = $normalStyleTemplate.Style
$ws.Range('E34').Value = '  +3.07%  '

$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.24'
$ws.Range('D36').Style = $normalStyleTemplate.Style
$ws.Range('E36').Value = '  +3.92%  '

$ws.Range('D37').Value = '3.710.73'
$ws.Range('E37').Value = '  -3.95%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.38'
$ws.Range('D38').Style = $normalStyleTemplate.Style
$ws.Range('E38').Value = '  +2.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.64'
$ws.Range('D39').Style = $normalStyleTemplate.Style
$ws.Range('E39').Value = '  +9.16%  '

$ws.Range('E40').Value = '  +4.07%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.66'
$ws.Range('D41').Style = $normalStyleTemplate.Style
$ws.Range('E41').Value = '  +3.48%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.20'
$ws.Range('D42').Style = $normalStyleTemplate.Style
$ws.Range('E42').Value = '  +7.56%  '

$ws.Range('D43').Value = '0.0₃0678'
$ws.Range('E43').Value = '  +1.61%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.35'
$ws.Range('D44').Style = $normalStyleTemplate.Style
$ws.Range('E44').Value = '  -1.15%  '

$ws.Range('E45').Value = '  +3.77%  '

$ws.Range('E46').Value = '  +1.09%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.67'
$ws.Range('D47').Style = $normalStyleTemplate.Style
$ws.Range('E47').Value = '  +5.24%  '

$ws.Range('E48').Value = '  +1.35%  '

$ws.Range('E49').Value = '  -0.32%  '

$ws.Range('E50').Value = '  -2.42%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '130.45'
$ws.Range('D51').Style = $normalStyleTemplate.Style
$ws.Range('E51').Value = '  +0.45%  '
